# Add a "LoadingUI" column (new column H) to the Scene sheet, inserted
# before the existing SoundList column, and populate it with the
# wallpaper path used while loading each scene.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert the new column (shifts old H..K -> I..L) ---------------------
[void]$ws.Columns("H").Insert()
$ws.Columns("H").ColumnWidth = 37.8

# --- header cell -----------------------------------------------------------
$ws.Cells.Item(1, 8).Value = "LoadingUI"
$ws.Cells.Item(1, 8).Font.Name = "宋体"
$ws.Cells.Item(1, 8).Font.Size = 11
$ws.Cells.Item(1, 8).Font.ColorIndex = 1
$ws.Cells.Item(1, 8).Font.Family = 3

# --- data rows ---------------------------------------------------------
# rows 2,3,6,7 -> the "Caster" wallpaper; rows 4,5 -> "forest" wallpaper
foreach ($r in 2,3,6,7) {
    $c = $ws.Cells.Item($r, 8)
    $c.NumberFormat = "@"
    $c.Value = "UI/ChronoBlade_Caster_wallpaper"
}

$h4 = $ws.Cells.Item(4, 8)
$h4.NumberFormat = "@"
$h4.Font.Name = "宋体"
$h4.Font.Size = 11
$h4.Font.ColorIndex = 1
$h4.Font.Family = 3
$h4.Value = "UI/ChronoBlade_forest_wallpaper"
# colour only the scene-specific suffix, matching the source file's
# rich-text run split ("UI/" stays default, the rest gets an explicit font)
$run = $h4.Characters(4, 28)
$run.Font.Name = "宋体"
$run.Font.Size = 11
$run.Font.ColorIndex = 1

# row 5 uses the identical rich-text wallpaper string - assign the same
# flattened text so it resolves back to the one shared-string entry
$h5 = $ws.Cells.Item(5, 8)
$h5.NumberFormat = "@"
$h5.Font.Name = "宋体"
$h5.Font.Size = 11
$h5.Font.ColorIndex = 1
$h5.Font.Family = 3
$h5.Value = $h4.Value2

# --- restore the active selection --------------------------------------
[void]$ws.Range("H1").Select()
